$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$ws.Activate()

# --- Remove the 15 rows in the log that never got a "what was done" entry ---
# (dates 45417,45418,45419,45421,45425,45426,45428,45432,45433,45434,45435,
#  45436,45437,45438,45439 had no matching B-column text). Deleting from the
# bottom up keeps the remaining row numbers stable while we iterate.
$emptyRows = @(25,24,23,22,21,20,19,18,14,12,11,7,5,4,3)
foreach ($r in $emptyRows) {
    $ws.Rows.Item($r).Delete()
}

# After the deletions the former row 55 (last log entry) is now row 40.
# Update its text to mention the new fog-effect work.
$ws.Range("B40").Value2 = "create a torch and candle - add optimisation for shading, add fog effect"

# --- Append the new log entry for the latest day of work ---
$ws.Range("A40:B40").Copy()
$ws.Range("A41:B41").PasteSpecial(-4122)
$ws.Range("A41").Value2 = 45512
$ws.Range("B41").Value2 = "work with textures, create new simple meshes in godot for walls, corridors, columns. Debug lights and shadows. "

# Restore the view roughly where the author left it
$ws.Range("B25").Select()
$excel.ActiveWindow.ScrollRow = 12

Write-Output "Logs sheet updated"
